$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.982.75"
$ws.Range("E2").Value = "  -7.85%  "
$ws.Range("D3").Value = "1.407.16"
$ws.Range("E3").Value = "  -8.52%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'273.31"
$ws.Range("E6").Value = "  -5.67%  "
$ws.Range("E7").Value = "  -6.61%  "
$ws.Range("D8").Value = "'0.3132"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").Value = "'39.70"
$ws.Range("E9").Value = "  -6.93%  "
$ws.Range("D10").Value = "'1.013"
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("D11").Value = "'0.06507"
$ws.Range("E11").Value = "  -9.53%  "
$ws.Range("D12").Value = "'0.9993"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'5.483"
$ws.Range("E13").Value = "  -4.69%  "
$ws.Range("D14").Value = "'17.35"
$ws.Range("E14").Value = "  -6.03%  "
$ws.Range("D15").Value = "'6.154"
$ws.Range("E15").Value = "  -7.52%  "
$ws.Range("D16").Value = "1.403.69"
$ws.Range("E16").Value = "  -8.90%  "
$ws.Range("D17").Value = "'0.00001014"
$ws.Range("E17").Value = "  -7.79%  "
$ws.Range("D18").Value = "'0.05708"
$ws.Range("E18").Value = "  -13.53%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'0.9999"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'70.96"
$ws.Range("E20").Value = "  -16.05%  "
$ws.Range("D21").Value = "'5.573"
$ws.Range("E21").Value = "  -9.57%  "
$ws.Range("D22").Value = "'14.72"
$ws.Range("E22").Value = "  -5.62%  "
$ws.Range("D23").Value = "'11.07"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").Value = "'2.254"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").Value = "19.975.11"
$ws.Range("E25").Value = "  -7.94%  "
$ws.Range("D26").Value = "'2.238"
$ws.Range("E26").Value = "  -6.92%  "
$ws.Range("D27").Value = "'135.80"
$ws.Range("E27").Value = "  -10.33%  "
$ws.Range("D28").Value = "'16.95"
$ws.Range("E28").Value = "  -8.33%  "
$ws.Range("D29").Value = "1.564.72"
$ws.Range("E29").Value = "  -8.78%  "
$ws.Range("D30").Value = "'109.47"
$ws.Range("E30").Value = "  -7.15%  "
$ws.Range("D31").Value = "'4.096"
$ws.Range("E31").Value = "  -15.81%  "
$ws.Range("D32").Value = "'5.313"
$ws.Range("E32").Value = "  -13.50%  "
$ws.Range("D33").Value = "'0.8191"
$ws.Range("E33").Value = "  -17.02%  "
$ws.Range("D34").Value = "'0.07671"
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("D35").Value = "'8.436"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").Value = "'1.444"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").Value = "'0.05786"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "'4.835"
$ws.Range("E38").Value = "  -7.21%  "
$ws.Range("D39").Value = "'1.000"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.02077"
$ws.Range("D41").Value = "'10.41"
$ws.Range("E41").Value = "  -8.10%  "
$ws.Range("D42").Value = "'0.1899"
$ws.Range("E42").Value = "  -7.65%  "
$ws.Range("D43").Value = "'1.103"
$ws.Range("E43").Value = "  -7.50%  "
$ws.Range("D44").Value = "'0.5275"
$ws.Range("E44").Value = "  -9.88%  "
$ws.Range("D45").Value = "'12.26"
$ws.Range("E45").Value = "  -7.07%  "
$ws.Range("D46").Value = "'3.511"
$ws.Range("E46").Value = "  -5.86%  "
$ws.Range("D47").Value = "'0.5123"
$ws.Range("E47").Value = "  -8.63%  "
$ws.Range("D48").Value = "'111.62"
$ws.Range("E48").Value = "  -4.69%  "
$ws.Range("D49").Value = "'1.764"
$ws.Range("E49").Value = "  -7.13%  "
$ws.Range("E50").Value = "  -11.49%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.38%  "
